$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used by PasteSpecial to copy formatting only
$xlPasteFormats = -4122

# --- Step 1: bump A1:B1 to the new bold "header" style (fontId4 bold, no wrap) ---
$ws.Range("A1:B1").Font.Bold = $true

# --- Step 2: build the new bold header style once on A29/B29 (creates new font+xf entries),
#     then propagate it to A54/B54 purely by copying formats (no new styles created) ---
$ws.Range("A29").Font.Bold = $true
$ws.Range("B29").Font.Bold = $true
$ws.Range("B29").WrapText = $true

$ws.Range("A29").Copy()
$ws.Range("A54").PasteSpecial($xlPasteFormats)
$ws.Range("B29").Copy()
$ws.Range("B54").PasteSpecial($xlPasteFormats)

# --- Step 3: propagate the pre-existing data styles (1,2,3) from their original source cells
#     onto every new cell that needs them, then set values afterwards ---
# Row 29
$ws.Cells.Item(29, 1).Value = "yield contract addresses"
$ws.Cells.Item(29, 2).Value = "name"

# Row 30
$ws.Range("B5").Copy()
$ws.Range("B30").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(30, 1).Value = "0x415747EE98D482e6dD9B431fa76Ad5553744F247"
$ws.Cells.Item(30, 2).Value = "`nDAI - eYyvDAI`n Apr 30, 2022"
$ws.Rows.Item(30).RowHeight = 43.2

# Row 31
$ws.Range("B3").Copy()
$ws.Range("B31").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(31, 1).Value = "0x8E9d636BbE6939BD0F52849afc02C0c66F6A3603"
$ws.Cells.Item(31, 2).Value = "LUSD3CRV-f - eYyvCurveLUSD`n Apr 30, 2022"
$ws.Rows.Item(31).RowHeight = 64.8

# Row 32
$ws.Range("B4").Copy()
$ws.Range("B32").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(32, 1).Value = "0xCF354603A9AEbD2Ff9f33E1B04246d8Ea204ae95"
$ws.Cells.Item(32, 2).Value = "`nWBTC - eYyvWBTC`n Apr 30, 2022"
$ws.Rows.Item(32).RowHeight = 34.2

# Row 33
$ws.Range("B3").Copy()
$ws.Range("B33").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(33, 1).Value = "0x7173b184525feAD2fFbde5FBe6FCB65Ea8246eE7"
$ws.Cells.Item(33, 2).Value = "USDC - eYyvUSDC`n Apr 29, 2022"
$ws.Rows.Item(33).RowHeight = 43.2

# Row 34
$ws.Range("B4").Copy()
$ws.Range("B34").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(34, 1).Value = "0x4aBB6FD289fA70056CFcB58ceBab8689921eB922"
$ws.Cells.Item(34, 2).Value = "crv3crypto - eYyvcrv3crypto`nApr 29, 2022"
$ws.Rows.Item(34).RowHeight = 22.8

# Row 35
$ws.Range("B5").Copy()
$ws.Range("B35").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(35, 1).Value = "0x7C9cF12d783821d5C63d8E9427aF5C44bAd92445"
$ws.Cells.Item(35, 2).Value = "USDC - eYyvUSDC`nDec 17, 2021"
$ws.Rows.Item(35).RowHeight = 28.8

# Row 36
$ws.Range("B5").Copy()
$ws.Range("B36").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(36, 1).Value = "0x062F38735AAC32320DB5e2DBBEb07968351D7C72"
$ws.Cells.Item(36, 2).Value = "steCRV - eYyvcrvSTETH`n Apr 15, 2022"
$ws.Rows.Item(36).RowHeight = 28.8

# Row 37
$ws.Range("B5").Copy()
$ws.Range("B37").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(37, 1).Value = "0xB70c25D96EF260eA07F650037Bf68F5d6583885e"
$ws.Cells.Item(37, 2).Value = "DAI - eYyvDAI`nJan 28, 2022"
$ws.Rows.Item(37).RowHeight = 28.8

# Row 38
$ws.Range("B3").Copy()
$ws.Range("B38").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(38, 1).Value = "0x4212bE3C7b255bA4B29705573ABD023cdcE21542"
$ws.Cells.Item(38, 2).Value = "steCRV - eYyvcrvSTETH`nJan 28, 2022"
$ws.Rows.Item(38).RowHeight = 64.8

# Row 39
$ws.Range("B4").Copy()
$ws.Range("B39").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(39, 1).Value = "0x9e030b67a8384cbba09D5927533Aa98010C87d91"
$ws.Cells.Item(39, 2).Value = "USDC - eYyvUSDC`n Jan 28, 2022"
$ws.Rows.Item(39).RowHeight = 22.8

# Row 40
$ws.Range("B5").Copy()
$ws.Range("B40").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(40, 1).Value = "0x7320d680Ca9BCE8048a286f00A79A2c9f8DCD7b3"
$ws.Cells.Item(40, 2).Value = "WBTC - eYyvWBTC`nNov 26, 2021"
$ws.Rows.Item(40).RowHeight = 28.8

# Row 41
$ws.Range("B5").Copy()
$ws.Range("B41").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(41, 1).Value = "0xd16847480D6bc218048CD31Ad98b63CC34e5c2bF"
$ws.Cells.Item(41, 2).Value = "crv3crypto - eYyvcrv3crypto`nNov 13, 2021"
$ws.Rows.Item(41).RowHeight = 28.8

# Row 42
$ws.Range("B5").Copy()
$ws.Range("B42").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(42, 1).Value = "0x2D6e3515C8b47192Ca3913770fa741d3C4Dac354"
$ws.Cells.Item(42, 2).Value = "USDC - eYyvUSDC`nOct 29, 2021"
$ws.Rows.Item(42).RowHeight = 28.8

# Row 43
$ws.Range("B5").Copy()
$ws.Range("B43").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(43, 1).Value = "0xE54B3F5c444a801e61BECDCa93e74CdC1C4C1F90"
$ws.Cells.Item(43, 2).Value = "DAI - eYyvDAI`n Oct 16, 2021"
$ws.Rows.Item(43).RowHeight = 28.8

# Row 44
$ws.Range("B5").Copy()
$ws.Range("B44").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(44, 1).Value = "0xD5D7bc115B32ad1449C6D0083E43C87be95F2809"
$ws.Cells.Item(44, 2).Value = "steCRV - eYyvcrvSTETH`nOct 16, 2021"
$ws.Rows.Item(44).RowHeight = 28.8

# Row 45
$ws.Range("B5").Copy()
$ws.Range("B45").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(45, 1).Value = "0xF94A7Df264A2ec8bCEef2cFE54d7cA3f6C6DFC7a"
$ws.Cells.Item(45, 2).Value = "crvTricrypto - eYyvCrvTriCrypto`nAug 15, 2021"
$ws.Rows.Item(45).RowHeight = 28.8

# Row 46
$ws.Range("B5").Copy()
$ws.Range("B46").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(46, 1).Value = "0x67F8FCb9D3c463da05DE1392EfDbB2A87F8599Ea"
$ws.Cells.Item(46, 2).Value = "LUSD3CRV-f - eYyvCurveLUSD`nDec 27, 2021"
$ws.Rows.Item(46).RowHeight = 28.8

# Row 47
$ws.Range("B5").Copy()
$ws.Range("B47").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(47, 1).Value = "0xDe620bb8BE43ee54d7aa73f8E99A7409Fe511084"
$ws.Cells.Item(47, 2).Value = "LUSD3CRV-f - eYyvCurveLUSD`nSep 28, 2021"
$ws.Rows.Item(47).RowHeight = 28.8

# Row 48
$ws.Range("B5").Copy()
$ws.Range("B48").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(48, 1).Value = "0x63E9B50DD3eB63BfBF93B26F57b9EFB574e59576"
$ws.Cells.Item(48, 2).Value = "alUSD3CRV-f - ePyvCurve-alUSD`n Apr 30, 2022"
$ws.Rows.Item(48).RowHeight = 28.8

# Row 49
$ws.Range("B5").Copy()
$ws.Range("B49").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(49, 1).Value = "0x6FE95FafE2F86158c77Bf18350672D360BfC78a2"
$ws.Cells.Item(49, 2).Value = "`nMIM-3LP3CRV-f - eYyvCurve-MIM`n Apr 29, 2022"
$ws.Rows.Item(49).RowHeight = 43.2

# Row 50
$ws.Range("B5").Copy()
$ws.Range("B50").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(50, 1).Value = "0x5fA3ce1fB47bC8A29B5C02e2e7167799BBAf5F41"
$ws.Cells.Item(50, 2).Value = "eursCRV - eYyvCurve-EURS`n Feb 12, 2022"
$ws.Rows.Item(50).RowHeight = 28.8

# Row 51
$ws.Range("B5").Copy()
$ws.Range("B51").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(51, 1).Value = "0x1D310a6238e11c8BE91D83193C88A99eB66279bE"
$ws.Cells.Item(51, 2).Value = "MIM-3LP3CRV-f - eYyvCurve-MIM`nFeb 11, 2022"
$ws.Rows.Item(51).RowHeight = 28.8

# Row 52
$ws.Range("B5").Copy()
$ws.Range("B52").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(52, 1).Value = "0x802d0f2f4b5f1fb5BfC9b2040a703c1464e1D4CB"
$ws.Cells.Item(52, 2).Value = "alUSD3CRV-f - eYyvCurve-alUSD`nJan 28, 2022"
$ws.Rows.Item(52).RowHeight = 28.8

# Row 54
$ws.Cells.Item(54, 1).Value = "Yield contract address"
$ws.Cells.Item(54, 2).Value = "name"

# Row 55
$ws.Range("B5").Copy()
$ws.Range("B55").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(55, 1).Value = "0xCFe60a1535ecc5B0bc628dC97111C8bb01637911"
$ws.Cells.Item(55, 2).Value = "USDC-sep17-2022"

# Row 56
$ws.Cells.Item(56, 1).Value = "0x52C9886d5D87B0f06EbACBEff750B5Ffad5d17d9"
$ws.Cells.Item(56, 2).Value = "USDC-apr29-2022"

# Row 57
$ws.Cells.Item(57, 1).Value = "0x2c72692E94E757679289aC85d3556b2c0f717E0E"
$ws.Cells.Item(57, 2).Value = "DAI-apr30-2022"

# Row 58
$ws.Cells.Item(58, 1).Value = "0x49e9e169f0B661Ea0A883f490564F4CC275123Ed"
$ws.Cells.Item(58, 2).Value = "WBTC-apr30-2022"

$excel.CutCopyMode = $false

# --- Step 4: restore view/selection state (best-effort; topLeftCell is not serialized by this engine) ---
$ws.Range("D63").Select()
